# "label ms on it"
# The CPU-time columns were measured/labeled in milliseconds, not seconds, and
# the three "caseN" chart titles are clarified to show that the x-axis / input
# size is expressed in ms. This script:
#   1. Relabels the "CPU time(s)" column headers to "CPU time(ms)".
#   2. Appends " (ms-inputsize)" to each of the six chart titles.
#   3. Reproduces the cosmetic view changes (zoom level, selected cell,
#      a few column widths) that accompanied the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relabel the CPU time header cells (row 4) from "CPU time(s)" to "CPU time(ms)"
$ws.Range("C4").Value = "CPU time(ms)"
$ws.Range("E4").Value = "CPU time(ms)"
$ws.Range("G4").Value = "CPU time(ms)"
$ws.Range("I4").Value = "CPU time(ms)"

# 2. Update all six chart titles to note the input size is in ms
$chartObjects = $ws.ChartObjects()
for ($i = 1; $i -le $chartObjects.Count; $i++) {
    $chart = $chartObjects.Item($i).Chart
    $title = $chart.ChartTitle
    $currentText = $title.Text
    if ($currentText -notlike "*(ms-inputsize)*") {
        $title.Text = "$currentText (ms-inputsize)"
    }
}

# 3. Cosmetic view tweaks that came along with the edit
$win = $excel.ActiveWindow
$win.Zoom = 89
$ws.Range("W35").Select()

# A handful of columns were widened slightly so the new, longer header text fits
$ws.Columns.Item(3).ColumnWidth = 10.830729166666666
$ws.Columns.Item(5).ColumnWidth = 10.830729166666666
$ws.Columns.Item(7).ColumnWidth = 10.998697916666666
$ws.Columns.Item(9).ColumnWidth = 10.830729166666666
